$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "09/01/2026 14:53:13"
$ws.Range("B2").Value = "14:53"
$ws.Range("C2").Value = "14:58"

# Row 3
$ws.Range("A3").Value = "09/01/2026 14:53:13"
$ws.Range("B3").Value = "14:53"
$ws.Range("C3").Value = "15:08"

# Row 4
$ws.Range("A4").Value = "09/01/2026 14:53:13"
$ws.Range("B4").Value = "14:53"
$ws.Range("C4").Value = "15:19"
